$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q header with the custom date format style
$ws.Range("Q1").Value = "FECHA_INGRESO_DUAL"
$ws.Range("Q1").NumberFormat = "dd/mm/yyyy;@"

# Column width for Q (bestFit-like autofit)
$ws.Columns.Item(17).AutoFit()

# Select Q2 to match the recorded UI selection state
$ws.Range("Q2").Select()
